# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update column G ("K") values for each row with the newly-recomputed counts,
# and (for row 37, where the recompute also affected IP / IF) update H and J too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = @{ G = 0 }
    3  = @{ G = 0 }
    4  = @{ G = 0 }
    5  = @{ G = 1 }
    6  = @{ G = 0 }
    7  = @{ G = 2 }
    8  = @{ G = 1 }
    9  = @{ G = 0 }
    10 = @{ G = 1 }
    11 = @{ G = 2 }
    12 = @{ G = 1 }
    13 = @{ G = 3 }
    14 = @{ G = 1 }
    15 = @{ G = 1 }
    16 = @{ G = 2 }
    17 = @{ G = 2 }
    18 = @{ G = 2 }
    19 = @{ G = 1 }
    20 = @{ G = 0 }
    21 = @{ G = 1 }
    22 = @{ G = 1 }
    23 = @{ G = 1 }
    24 = @{ G = 0 }
    25 = @{ G = 3 }
    26 = @{ G = 2 }
    27 = @{ G = 1 }
    28 = @{ G = 0 }
    29 = @{ G = 4 }
    30 = @{ G = 2 }
    31 = @{ G = 1 }
    32 = @{ G = 0 }
    33 = @{ G = 0 }
    34 = @{ G = 1 }
    35 = @{ G = 2 }
    36 = @{ G = 0 }
    37 = @{ G = 0; H = 7; J = 11 }
    38 = @{ G = 1 }
    39 = @{ G = 0 }
    40 = @{ G = 1 }
    41 = @{ G = 1 }
    42 = @{ G = 1 }
    44 = @{ G = 1 }
    46 = @{ G = 0 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
